$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.232.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.140.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.84%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.15%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.139.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.20%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.38%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.66%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.86%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.657.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.23%  "

# Row 16
$ws.Range("E16").Value = "  -2.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.134.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.234.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.00%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.96%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.81%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.89%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("E27").Value = "  +0.05%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.72%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.81%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.28%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.20%  "

# Row 33
$ws.Range("E33").Value = "  -4.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.50%  "

# Row 35
$ws.Range("E35").Value = "  -7.40%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.26%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0692"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.38%  "

# Row 39
$ws.Range("E39").Value = "  -4.55%  "

# Row 40
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.81%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.23%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "391.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.78%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.111"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.25%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.747.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.39%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.250"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.55%  "

# Row 46
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.48%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.25%  "

# Row 51
$ws.Range("E51").Value = "  -3.40%  "
